# The "Name" column (A) is being removed; the remaining "Ticker" and
# "Actual" columns shift left to become columns A and B respectively.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("A").Delete()
